$wb = $excel.ActiveWorkbook

# --- Sheet: inputdisp ---
$ws1 = $wb.Worksheets.Item("inputdisp")
$ws1.Range("G2").Value = 20
$ws1.Range("G3").Value = 6
$ws1.Activate()
$ws1.Range("G3").Select()

# --- Sheet: endofpipe ---
$ws2 = $wb.Worksheets.Item("endofpipe")
$ws2.Range("G2").Value = 10
$ws2.Activate()
$ws2.Range("G3").Select()

# --- Sheet: inputprices ---
$ws3 = $wb.Worksheets.Item("inputprices")
$ws3.Range("B2").Value = 3
$ws3.Range("B3").Value = 0.1
$ws3.Range("B4").Value = 3
$ws3.Range("B5").Value = 3
$ws3.Activate()
$ws3.Range("E7").Select()

# Restore the originally active sheet (inputdisp, tabSelected in sheet1.xml)
$ws1.Activate()
$ws1.Range("G3").Select()
